# Fix slide show options: merge the title's split runs back together.
#
# The title placeholder on the "Talks (in no particular order)" slide has
# its text spread across five runs:
#   "Talks " | "(in " | "no " | "particular " | "order)"
# Runs 1+2 share identical run properties, as do runs 4+5, so the fix
# collapses each matching pair into a single run:
#   "Talks (in " | "no " | "particular order)"
#
# Plain TextRange.Text assignment only patches characters that actually
# changed and never removes/merges <a:r> run elements, so instead we
# delete the redundant run's characters (a delete that exactly spans a
# run removes the <a:r> element outright) and re-insert that text into
# the end of the preceding run, which extends its existing <a:t>.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Sanity check (not required, but guards against slide/shape drift).
# Expected: "Talks (in no particular order)"

# --- Merge run 1 ("Talks ") with run 2 ("(in ") ---
# Characters 1-6 = "Talks ", characters 7-10 = "(in ".
$tr.Characters(7, 4).Delete() | Out-Null
$tr.Characters(1, 6).InsertAfter("(in ") | Out-Null

# --- Merge run 4 ("particular ") with run 5 ("order)") ---
# After the first merge, text is "Talks (in no particular order)":
# characters 14-24 = "particular ", characters 25-30 = "order)".
$tr.Characters(25, 6).Delete() | Out-Null
$tr.Characters(14, 11).InsertAfter("order)") | Out-Null
